$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D (Price) holds text-formatted values even though many of them
# look numeric (e.g. "287.01"). Force the whole range to Text format
# before writing so Excel does not silently coerce them into floating
# point numbers (which would also drop trailing/leading zeros).
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = "22.405.02"
$ws.Range("E2").Value = "  -0.32%  "

$ws.Range("D3").Value = "1.564.00"
$ws.Range("E3").Value = "  -0.52%  "

$ws.Range("E4").Value = "  +0.00%  "

$ws.Range("E5").Value = "  +0.00%  "

$ws.Range("D6").Value = "287.01"
$ws.Range("E6").Value = "  -1.47%  "

$ws.Range("D7").Value = "0.3660"
$ws.Range("E7").Value = "  -2.47%  "

$ws.Range("D8").Value = "49.67"
$ws.Range("E8").Value = "  -0.41%  "

$ws.Range("D9").Value = "0.3356"
$ws.Range("E9").Value = "  -1.52%  "

$ws.Range("E10").Value = "  -1.90%  "

$ws.Range("D11").Value = "0.07416"
$ws.Range("E11").Value = "  -2.15%  "

$ws.Range("E12").Value = "  -0.02%  "

$ws.Range("D13").Value = "20.92"
$ws.Range("E13").Value = "  -2.36%  "

$ws.Range("D14").Value = "5.941"
$ws.Range("E14").Value = "  -1.30%  "

$ws.Range("D15").Value = "6.885"
$ws.Range("E15").Value = "  -1.09%  "

$ws.Range("D16").Value = "1.564.44"
$ws.Range("E16").Value = "  -0.56%  "

$ws.Range("D17").Value = "0.00001102"
$ws.Range("E17").Value = "  -1.90%  "

$ws.Range("D18").Value = "89.09"
$ws.Range("E18").Value = "  -2.25%  "

$ws.Range("D19").Value = "0.06746"
$ws.Range("E19").Value = "  +0.04%  "

$ws.Range("E20").Value = "  -0.04%  "

$ws.Range("D21").Value = "6.307"
$ws.Range("E21").Value = "  +0.52%  "

$ws.Range("D22").Value = "16.08"
$ws.Range("E22").Value = "  -2.13%  "

$ws.Range("E23").Value = "  -1.98%  "

$ws.Range("D24").Value = "22.393.10"
$ws.Range("E24").Value = "  -0.38%  "

$ws.Range("D25").Value = "2.373"
$ws.Range("E25").Value = "  +1.81%  "

$ws.Range("D26").Value = "2.543"
$ws.Range("E26").Value = "  -2.26%  "

$ws.Range("D27").Value = "149.72"
$ws.Range("E27").Value = "  +0.88%  "

$ws.Range("E28").Value = "  -2.85%  "

$ws.Range("D29").Value = "5.007"
$ws.Range("E29").Value = "  +0.19%  "

$ws.Range("D30").Value = "123.24"
$ws.Range("E30").Value = "  -2.33%  "

$ws.Range("D31").Value = "1.740.85"
$ws.Range("E31").Value = "  -0.48%  "

$ws.Range("D32").Value = "1.073"
$ws.Range("E32").Value = "  +2.08%  "

$ws.Range("D33").Value = "6.129"
$ws.Range("E33").Value = "  -0.60%  "

$ws.Range("D34").Value = "1.996"
$ws.Range("E34").Value = "  +0.91%  "

$ws.Range("D35").Value = "9.624"
$ws.Range("E35").Value = "  -2.54%  "

$ws.Range("D36").Value = "0.08281"
$ws.Range("E36").Value = "  -2.16%  "

$ws.Range("D37").Value = "0.02400"
$ws.Range("E37").Value = "  -3.03%  "

$ws.Range("D38").Value = "1.312"
$ws.Range("E38").Value = "  -4.90%  "


# Rows 39 and 40 swapped places (Algorand <-> Hedera), with refreshed
# price/volume figures for each coin in its new row position.
$ws.Range("B39").Value = "Hedera"
$ws.Range("C39").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D39").Value = "0.06410"
$ws.Range("E39").Value = "  -2.39%  "

$ws.Range("B40").Value = "Algorand"
$ws.Range("C40").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("D40").Value = "0.2221"
$ws.Range("E40").Value = "  -3.22%  "

$ws.Range("D41").Value = "5.341"
$ws.Range("E41").Value = "  -2.76%  "

$ws.Range("D42").Value = "11.21"
$ws.Range("E42").Value = "  -1.86%  "

$ws.Range("D43").Value = "0.6083"
$ws.Range("E43").Value = "  -3.57%  "

$ws.Range("E44").Value = "  -0.01%  "

$ws.Range("D45").Value = "13.76"
$ws.Range("E45").Value = "  -2.81%  "

$ws.Range("D46").Value = "3.771"
$ws.Range("E46").Value = "  -1.21%  "

$ws.Range("D47").Value = "0.5742"
$ws.Range("E47").Value = "  -2.61%  "

$ws.Range("D48").Value = "2.019"
$ws.Range("E48").Value = "  -3.97%  "

$ws.Range("D49").Value = "124.95"
$ws.Range("E49").Value = "  -3.83%  "

$ws.Range("E50").Value = "  -1.06%  "

$ws.Range("D51").Value = "0.07236"
$ws.Range("E51").Value = "  -1.41%  "

